$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update the Contact value (row 10)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for "Jurisdiction" right after "Contact" (row 11), shifting
# the remaining metadata rows down by one.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
# Leading apostrophe forces Excel to store an explicit (empty) text value
# for the Jurisdiction value cell instead of leaving it a blank cell.
$ws.Range("B11").Value = "'"

# The freshly inserted row does not carry the same cell style as the
# surrounding data rows, so copy the formatting from the row above (Contact)
# onto it.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
